$wb = $excel.ActiveWorkbook

# --- 1. "opmaak" sheet: insert a new row with header_template / Totaal [naam] [jaar] ---
$opmaak = $wb.Worksheets.Item("opmaak")
$opmaak.Activate()
$opmaak.Rows.Item(22).Insert() | Out-Null
$opmaak.Cells.Item(22, 1).Value = "header_template"
$opmaak.Cells.Item(22, 2).Value = "Totaal [naam] [jaar]"
$opmaak.Range("B23").Select() | Out-Null

# --- 2. Add new "labelcorrectie" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$labelcorrectie = $wb.Worksheets.Add($null, $lastSheet)
$labelcorrectie.Name = "labelcorrectie"

# Column widths (values chosen so the engine's char-width -> stored-width
# conversion lands on the target stored widths of 12 / 21 / 22.7109375 /
# 18.7109375 / 16)
$labelcorrectie.Columns.Item(1).ColumnWidth = 11.1666666666667
$labelcorrectie.Columns.Item(2).ColumnWidth = 20.1666666666667
$labelcorrectie.Columns.Item(3).ColumnWidth = 21.8333333333333
$labelcorrectie.Columns.Item(4).ColumnWidth = 17.8333333333333
$labelcorrectie.Columns.Item(5).ColumnWidth = 15.1666666666667

# Row 1 - header row (write left to right so shared strings are minted in this order)
$labelcorrectie.Cells.Item(1, 1).Value = "var"
$labelcorrectie.Cells.Item(1, 2).Value = "var_label"
$labelcorrectie.Cells.Item(1, 3).Value = "antwoord_waarde"
$labelcorrectie.Cells.Item(1, 4).Value = "antwoord_oud"
$labelcorrectie.Cells.Item(1, 5).Value = "antwoord_nieuw"

# Row 3 filled next (example/test row)
$labelcorrectie.Cells.Item(3, 1).Value = "AGHHA401"
$labelcorrectie.Cells.Item(3, 2).Value = "Test"

# Row 2
$labelcorrectie.Cells.Item(2, 1).Value = "dagbesteding"
$labelcorrectie.Cells.Item(2, 2).Value = "Abc"
$labelcorrectie.Cells.Item(2, 3).Value = 1
$labelcorrectie.Cells.Item(2, 5).Value = "Alleen opleiding"

# Row 4
$labelcorrectie.Cells.Item(4, 1).Value = "dagbesteding"
$labelcorrectie.Cells.Item(4, 4).Value = "Werkt, volgt geen opleiding"
$labelcorrectie.Cells.Item(4, 5).Value = "Geen opleiding"

# Formatting for A3 (wrap text, white fill, number format "0")
$labelcorrectie.Cells.Item(3, 1).Interior.Color = 16777215
$labelcorrectie.Cells.Item(3, 1).NumberFormat = "0"
$labelcorrectie.Cells.Item(3, 1).WrapText = $true

$labelcorrectie.Range("C8").Select()
